$wb = $excel.ActiveWorkbook

# --- "About" sheet: clear explicit (redundant, non-bold) font styling on A11 ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A11").Style = "Normal"

# --- "QSfHO" sheet: update the Quantization Size value and leave the cursor on B3 ---
$wsQSfHO = $wb.Worksheets.Item("QSfHO")
$wsQSfHO.Range("B2").Value = 1

# Reflect the final cell selection / active sheet as last edited by the author:
# cursor left on B3 of QSfHO, but the About sheet is the one shown/active on save.
$wsQSfHO.Activate()
$wsQSfHO.Range("B3").Select()
$wsAbout.Activate()
